# Rename header columns to match the new cases/products/technologies
# settings-file naming convention (c_Name/p_Name/t_Name -> cases_Name/
# products_Name/technologies_Name), then move the active selection
# from E11 to F11 to match the saved sheet view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "cases_Name"
$ws.Range("C1").Value = "products_Name"
$ws.Range("D1").Value = "technologies_Name"

$ws.Range("F11").Select()
